# Applies the edits described by the diff to the active workbook:
#  1. Adds a new header cell BD1 = "Odd_CS_4-4" (copying the header style from BC1)
#  2. Updates several odds values in row 2
#  3. Adds an (empty) placeholder cell BD2
#  4. Adds a whole new row 3 with match data (Avai vs Ponte Preta) including BD3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header cell BD1
# ---------------------------------------------------------------------------
$ws.Range("BD1").Value = "Odd_CS_4-4"
# Match the formatting used by the rest of the header row (bold, centered, bordered)
$ws.Range("BC1").Copy() | Out-Null
$ws.Range("BD1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Updated odds in row 2
# ---------------------------------------------------------------------------
$row2Updates = @{
    "G2"  = 1.53
    "H2"  = 3.7
    "I2"  = 7
    "L2"  = 7
    "Z2"  = 10
    "AG2" = 15
    "AH2" = 34
    "AJ2" = 81
    "AM2" = 3.25
    "AT2" = 10
    "AU2" = 81
    "AW2" = 41
}
foreach ($addr in $row2Updates.Keys) {
    $ws.Range($addr).Value = $row2Updates[$addr]
}

# ---------------------------------------------------------------------------
# 3. Empty placeholder cell BD2
# ---------------------------------------------------------------------------
$ws.Range("BD2").Font.Bold = $false

# ---------------------------------------------------------------------------
# 4. New row 3 (Avai vs Ponte Preta)
# ---------------------------------------------------------------------------
$row3Values = New-Object 'object[,]' 1,56
$row3Values[0,0] = "xrkDKQao"
$row3Values[0,1] = "22/11/2024"
$row3Values[0,2] = "20:00"
$row3Values[0,3] = "BRAZIL - SERIE B"
$row3Values[0,4] = "Avai"
$row3Values[0,5] = "Ponte Preta"
$row3Values[0,6] = 1.57
$row3Values[0,7] = 3.6
$row3Values[0,8] = 6.5
$row3Values[0,9] = 2.2
$row3Values[0,10] = 2.2
$row3Values[0,11] = 6.5
$row3Values[0,12] = 1.07
$row3Values[0,13] = 8.5
$row3Values[0,14] = 1.33
$row3Values[0,15] = 3.25
$row3Values[0,16] = 2.08
$row3Values[0,17] = 1.73
$row3Values[0,18] = 1.44
$row3Values[0,19] = 2.63
$row3Values[0,20] = 2.1
$row3Values[0,21] = 1.67
$row3Values[0,22] = 6
$row3Values[0,23] = 6.5
$row3Values[0,24] = 9
$row3Values[0,25] = 11
$row3Values[0,26] = 15
$row3Values[0,27] = 34
$row3Values[0,28] = 8.5
$row3Values[0,29] = 7
$row3Values[0,30] = 19
$row3Values[0,31] = 67
$row3Values[0,32] = 13
$row3Values[0,33] = 29
$row3Values[0,34] = 21
$row3Values[0,35] = 67
$row3Values[0,36] = 51
$row3Values[0,37] = 51
$row3Values[0,38] = 3.4
$row3Values[0,39] = 8
$row3Values[0,40] = 21
$row3Values[0,41] = 26
$row3Values[0,42] = 51
$row3Values[0,43] = 201
$row3Values[0,44] = 2.63
$row3Values[0,45] = 9.5
$row3Values[0,46] = 67
$row3Values[0,47] = 7.5
$row3Values[0,48] = 34
$row3Values[0,49] = 41
$row3Values[0,50] = 126
$row3Values[0,51] = 151
$row3Values[0,52] = 81
$row3Values[0,53] = 81
$row3Values[0,54] = 351
$row3Values[0,55] = 900

$ws.Range("A3:BD3").Value = $row3Values
